$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.5925266903914591
$ws1.Range("C2").Value = 0.08641975308641975
$ws1.Range("D2").Value = 0.75
$ws1.Range("E2").Value = 0.1549815498154982
$ws1.Range("F2").Value = 0.2957746478873239
$ws1.Range("G2").Value = 0.5790031813361611
$ws1.Range("H2").Value = 0.7796615837346175
$ws1.Range("I2").Value = 21
$ws1.Range("J2").Value = 222
$ws1.Range("K2").Value = 312
$ws1.Range("L2").Value = 7

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")

# row 2 (label 0)
$ws2.Range("B2").Value = 0.9780564263322884
$ws2.Range("C2").Value = 0.5842696629213483
$ws2.Range("D2").Value = 0.731535756154748

# row 3 (label 1)
$ws2.Range("B3").Value = 0.08641975308641975
$ws2.Range("C3").Value = 0.75
$ws2.Range("D3").Value = 0.1549815498154982

# row 4 (accuracy)
$ws2.Range("B4").Value = 0.5925266903914591
$ws2.Range("C4").Value = 0.5925266903914591
$ws2.Range("D4").Value = 0.5925266903914591
$ws2.Range("E4").Value = 0.5925266903914591

# row 5 (macro avg)
$ws2.Range("B5").Value = 0.532238089709354
$ws2.Range("C5").Value = 0.6671348314606742
$ws2.Range("D5").Value = 0.4432586529851231

# row 6 (weighted avg)
$ws2.Range("B6").Value = 0.9336332468823163
$ws2.Range("C6").Value = 0.5925266903914591
$ws2.Range("D6").Value = 0.702810635554216

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

# row 2 (Actual 0)
$ws3.Range("B2").Value = 312
$ws3.Range("C2").Value = 222

# row 3 (Actual 1)
$ws3.Range("B3").Value = 7
$ws3.Range("C3").Value = 21
